$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "data last refreshed" timestamp in A1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 2 de Octubre de 2020 a las 00:11"

# Swap the two country rows whose label position changed in the shared-string
# table (Santa Lucia now sorts before Nueva Caledonia; Montserrat now sorts
# before Islas Malvinas) by swapping the two full rows' contents.
function Swap-Rows($ws, $r1, $r2) {
    for ($c = 1; $c -le 8; $c++) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}

Swap-Rows $ws 207 208
Swap-Rows $ws 215 216

# Updated COVID-19 stats (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes) for the countries whose
# figures changed.
function Set-RowStats($ws, $r, $vals) {
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 2).Value = $vals[$i]
    }
}

Set-RowStats $ws 4   @(7489214, 41932, 4726829, 2549888, 0, 757, 212497)
Set-RowStats $ws 8   @(835339, 5660, 751691, 57452, 0, 198, 26196)
Set-RowStats $ws 29  @(160453, 1695, 136246, 14889, 0, 21, 9318)
Set-RowStats $ws 57  @(71374, 510, 65550, 5570, 0, 3, 254)
Set-RowStats $ws 75  @(35208, 683, 19894, 14755, 0, 8, 559)
Set-RowStats $ws 134 @(4247, 47, 1117, 2928, 0, 2, 202)
Set-RowStats $ws 161 @(1809, 25, 1353, 408, 0, 0, 48)
